$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.687.80'
$ws.Range("E2").Value = '  -2.19%  '
$ws.Range("D3").Value = '2.302.98'
$ws.Range("E3").Value = '  -4.43%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '546.42'
$ws.Range("E5").Value = '  -1.42%  '
$ws.Range("D6").Value = '131.87'
$ws.Range("E6").Value = '  -3.07%  '
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("E8").Value = '  -2.38%  '
$ws.Range("D9").Value = '2.302.06'
$ws.Range("E9").Value = '  -4.42%  '
$ws.Range("E10").Value = '  -2.54%  '
$ws.Range("E13").Value = '  -4.99%  '
$ws.Range("D14").Value = '23.98'
$ws.Range("E14").Value = '  -2.75%  '
$ws.Range("D15").Value = '2.710.96'
$ws.Range("E15").Value = '  -4.52%  '
$ws.Range("D16").Value = '58.689.12'
$ws.Range("E17").Value = '  -3.09%  '
$ws.Range("D18").Value = '2.303.63'
$ws.Range("E18").Value = '  -4.23%  '
$ws.Range("E19").Value = '  -4.39%  '
$ws.Range("E20").Value = '  -4.47%  '
$ws.Range("D21").Value = '314.56'
$ws.Range("E21").Value = '  -3.97%  '
$ws.Range("E22").Value = '  -4.03%  '
$ws.Range("E23").Value = '  +0.12%  '
$ws.Range("D24").Value = '63.34'
$ws.Range("E24").Value = '  -2.19%  '
$ws.Range("E25").Value = '  -5.97%  '
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  +0.02%  '
$ws.Range("E27").Value = '  -5.73%  '
$ws.Range("E28").Value = '  -5.76%  '
$ws.Range("E29").Value = '  -2.46%  '
$ws.Range("D30").Value = '169.29'
$ws.Range("E30").Value = '  -0.85%  '
$ws.Range("D31").Value = '0.0₃0725'
$ws.Range("E31").Value = '  -5.62%  '
$ws.Range("E32").Value = '  -0.35%  '
$ws.Range("E33").Value = '  -5.54%  '
$ws.Range("D34").Value = '0.380'
$ws.Range("E34").Value = '  -5.06%  '
$ws.Range("D36").Value = '17.78'
$ws.Range("E36").Value = '  -3.52%  '
$ws.Range("E37").Value = '  +0.07%  '
$ws.Range("E38").Value = '  -4.86%  '
$ws.Range("E39").Value = '  -5.46%  '
$ws.Range("D40").Value = '38.06'
$ws.Range("E40").Value = '  -2.33%  '
$ws.Range("E41").Value = '  -5.28%  '
$ws.Range("D42").Value = '297.97'
$ws.Range("E42").Value = '  -7.52%  '
$ws.Range("D43").Value = '140.48'
$ws.Range("E43").Value = '  -3.92%  '
$ws.Range("E44").Value = '  -4.55%  '
$ws.Range("D45").Value = '0.0953'
$ws.Range("E45").Value = '  -1.08%  '
$ws.Range("E46").Value = '  -2.54%  '
$ws.Range("D47").Value = '0.556'
$ws.Range("E47").Value = '  -3.44%  '
$ws.Range("D48").Value = '18.45'
$ws.Range("E48").Value = '  -7.29%  '
$ws.Range("D49").Value = '0.0215'
$ws.Range("E49").Value = '  -2.71%  '
$ws.Range("D50").Value = '16.63'
$ws.Range("E50").Value = '  -4.05%  '
$ws.Range("E51").Value = '  -0.34%  '
